$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 48-49; this pushes the existing rows 48-65 down to 50-67,
# which already contain the correct (shifted) target values for those rows.
$ws.Rows("48:49").Insert()

# Row 48 (new record)
$ws.Cells.Item(48, 1).Value2 = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value2 = 44551
$ws.Cells.Item(48, 5).Value2 = 8
$ws.Cells.Item(48, 6).Value2 = 100112028
$ws.Cells.Item(48, 7).Value = "Sandia"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value2 = 500
$ws.Cells.Item(48, 11).Value2 = 2000
$ws.Cells.Item(48, 12).Value2 = 2200
$ws.Cells.Item(48, 13).Value2 = 2120
$ws.Cells.Item(48, 14).Value = "`$/unidad"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value2 = 2120
$ws.Cells.Item(48, 17).Value2 = 1
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# Row 49 (new record)
$ws.Cells.Item(49, 1).Value2 = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value2 = 44551
$ws.Cells.Item(49, 5).Value2 = 8
$ws.Cells.Item(49, 6).Value2 = 100112028
$ws.Cells.Item(49, 7).Value = "Sandia"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Segunda"
$ws.Cells.Item(49, 10).Value2 = 300
$ws.Cells.Item(49, 11).Value2 = 1800
$ws.Cells.Item(49, 12).Value2 = 1800
$ws.Cells.Item(49, 13).Value2 = 1800
$ws.Cells.Item(49, 14).Value = "`$/unidad"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value2 = 1800
$ws.Cells.Item(49, 17).Value2 = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"
